$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the duration of the value writes, then restore the
# default "Normal" style so the cells keep their original (unstyled) look -
# only the stored value/type changes, matching how the source file was built.
$ws.Range("A2:H2").NumberFormat = "@"
$ws.Range("A2").Value = "39.1"
$ws.Range("B2").Value = "210.5"
$ws.Range("C2").Value = "Gentoo"
$ws.Range("D2").Value = "0.45"
$ws.Range("E2").Value = "0.0"
$ws.Range("F2").Value = "0.55"
$ws.Range("G2").Value = "v1.0"
$ws.Range("H2").Value = "2025-05-04 20:35:00"
$ws.Range("A2:H2").Style = "Normal"

# Row 3 (the old Chinstrap prediction) is removed entirely.
$ws.Rows("3:3").Delete()
